$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ChromeMobile"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "SafariTablet"

$ws.Range("A4:B5").Select()
